# B6-PowerPoint.pptx edit
#  1) Re-colour the deck's (slide-master) theme from the "Integral" / Red
#     Violet palette to the stock Office-theme palette - dk1/lt1 (black/white)
#     are already identical between the two palettes, so only dk2, lt2 and
#     the six accents + hyperlink colours need touching.
#  2) Point the three tables in the deck (on the "B6" slides) at the
#     built-in "Office" table style instead of the custom Table_0 style.

$p = $ppt.ActivePresentation

# --- 1. Theme re-colour -----------------------------------------------
# Office theme colours (RRGGBB), in clrScheme slot order.
$officeTheme = @{
    1  = 0x000000   # dk1
    2  = 0xFFFFFF   # lt1
    3  = 0x44546A   # dk2
    4  = 0xE7E6E6   # lt2
    5  = 0x5B9BD5   # accent1
    6  = 0xED7D31   # accent2
    7  = 0xA5A5A5   # accent3
    8  = 0xFFC000   # accent4
    9  = 0x4472C4   # accent5
    10 = 0x70AD47   # accent6
    11 = 0x0563C1   # hlink
    12 = 0x954F72   # folHlink
}

# The deck's theme colours are reached through any slide's ThemeColorScheme
# (they all resolve back to the single shared slide-master theme part).
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

foreach ($idx in $officeTheme.Keys) {
    $rrggbb = $officeTheme[$idx]
    $r = ($rrggbb -band 0xFF0000) -shr 16
    $g = ($rrggbb -band 0x00FF00) -shr 8
    $b = ($rrggbb -band 0x0000FF)
    # PowerPoint's RGB long is packed &H00BBGGRR (standard VBA RGB()).
    $bgrLong = ($b * 65536) + ($g * 256) + $r
    $themeColors.Item($idx).RGB = $bgrLong
}

# --- 2. Table style -----------------------------------------------------
$newTableStyleId = "{3CF904C3-D55D-42B5-A540-DE9095189B05}"

foreach ($slide in $p.Slides) {
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}
